$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F5: test data email was rotated to a new value.
$ws.Range("F5").Value = "meeyaaken2@gmail.com"

# H5: numeric id bumped to a new value, and its border was brought in
# line with the rest of row 5 (which all carry a thin box border).
$ws.Range("H5").Value = 4373388
$ws.Range("H5").Borders.LineStyle = 1

# Add a reviewer note on F5 explaining why the email should track the
# PIDIntegration value.
$comment = $ws.Range("F5").AddComment("Nitin Sharma:" + [char]10 + "change email id according to PIDIntegration value")
